$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 236.7946826666667
$ws.Range("H2").Value = 710.384048
$ws.Range("I2").Value = 0.7123899543147419
$ws.Range("J2").Value = 0.7240508783182559
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 3291.29185679669
$ws.Range("R2").Value = 29621.62671117021
$ws.Range("S2").Value = 0.03457590024043512
$ws.Range("T2").Value = 0.03619728642602047
$ws.Range("G3").Value = 236.7946826666667
$ws.Range("H3").Value = 710.384048
$ws.Range("I3").Value = 0.7123899543147419
$ws.Range("J3").Value = 0.7240508783182559
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 16661.43073458874
$ws.Range("R3").Value = 149952.8766112986
$ws.Range("S3").Value = 0.175032781049914
$ws.Range("T3").Value = 0.18324068688159
$ws.Range("G4").Value = 236.7946826666667
$ws.Range("H4").Value = 710.384048
$ws.Range("I4").Value = 0.7123899543147419
$ws.Range("J4").Value = 0.7240508783182559
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 19418.97900990931
$ws.Range("R4").Value = 174770.8110891838
$ws.Range("S4").Value = 0.204001562374723
$ws.Range("T4").Value = 0.21356791676527
$ws.Range("G5").Value = 236.7946826666667
$ws.Range("H5").Value = 710.384048
$ws.Range("I5").Value = 0.7123899543147419
$ws.Range("J5").Value = 0.7240508783182559
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 5931.741491221011
$ws.Range("R5").Value = 35590.44894732606
$ws.Range("S5").Value = 0.06231452905915196
$ws.Range("T5").Value = 0.04349111840273159
$ws.Range("G6").Value = 236.7946826666667
$ws.Range("H6").Value = 710.384048
$ws.Range("I6").Value = 0.7123899543147419
$ws.Range("J6").Value = 0.7240508783182559
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 22509.2021081973
$ws.Range("R6").Value = 202582.8189737757
$ws.Range("S6").Value = 0.2364651815905178
$ws.Range("T6").Value = 0.2475538698426439
$ws.Range("I7").Value = 0.2358656137148928
$ws.Range("J7").Value = 0.2397264359793184
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 1089.715778579327
$ws.Range("R7").Value = 9807.442007213942
$ws.Range("S7").Value = 0.0114477553768986
$ws.Range("T7").Value = 0.01198458109351023
$ws.Range("I8").Value = 0.2358656137148928
$ws.Range("J8").Value = 0.2397264359793184
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("R8").Value = 49647.98711158023
$ws.Range("S8").Value = 0.05795170759008567
$ws.Range("T8").Value = 0.06066926801408761
$ws.Range("I9").Value = 0.2358656137148928
$ws.Range("J9").Value = 0.2397264359793184
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 6429.441311107068
$ws.Range("R9").Value = 57864.97179996361
$ws.Range("S9").Value = 0.06754299863000651
$ws.Range("T9").Value = 0.07071032859539179
$ws.Range("I10").Value = 0.2358656137148928
$ws.Range("J10").Value = 0.2397264359793184
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 1963.943818622127
$ws.Range("R10").Value = 11783.66291173276
$ws.Range("S10").Value = 0.0206317544918632
$ws.Range("T10").Value = 0.01439950026116626
$ws.Range("I11").Value = 0.2358656137148928
$ws.Range("J11").Value = 0.2397264359793184
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 7452.585114832869
$ws.Range("R11").Value = 67073.26603349582
$ws.Range("S11").Value = 0.07829139762603875
$ws.Range("T11").Value = 0.08196275801516252
$ws.Range("G12").Value = 0.6305213333333334
$ws.Range("H12").Value = 1.891564
$ws.Range("I12").Value = 0.001896905195629352
$ws.Range("J12").Value = 0.001927955138422806
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 8.763835853771557
$ws.Range("R12").Value = 78.87452268394401
$ws.Range("S12").Value = 0.00009206643694566525
$ws.Range("T12").Value = 0.00009638375762225587
$ws.Range("G13").Value = 0.6305213333333334
$ws.Range("H13").Value = 1.891564
$ws.Range("I13").Value = 0.001896905195629352
$ws.Range("J13").Value = 0.001927955138422806
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 44.36496378933556
$ws.Range("R13").Value = 399.28467410402
$ws.Range("S13").Value = 0.0004660657969249607
$ws.Range("T13").Value = 0.0004879212696517192
$ws.Range("G14").Value = 0.6305213333333334
$ws.Range("H14").Value = 1.891564
$ws.Range("I14").Value = 0.001896905195629352
$ws.Range("J14").Value = 0.001927955138422806
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 51.70758227935334
$ws.Range("R14").Value = 465.36824051418
$ws.Range("S14").Value = 0.000543201965779193
$ws.Range("T14").Value = 0.000568674626134315
$ws.Range("G15").Value = 0.6305213333333334
$ws.Range("H15").Value = 1.891564
$ws.Range("I15").Value = 0.001896905195629352
$ws.Range("J15").Value = 0.001927955138422806
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 15.79465177137534
$ws.Range("R15").Value = 94.767910628252
$ws.Range("S15").Value = 0.0001659270364770997
$ws.Range("T15").Value = 0.0001158052945050711
$ws.Range("G16").Value = 0.6305213333333334
$ws.Range("H16").Value = 1.891564
$ws.Range("I16").Value = 0.001896905195629352
$ws.Range("J16").Value = 0.001927955138422806
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 59.93602544491556
$ws.Range("R16").Value = 539.42422900424
$ws.Range("S16").Value = 0.0006296439595024328
$ws.Range("T16").Value = 0.0006591701905094451
$ws.Range("G17").Value = 16.059769
$ws.Range("H17").Value = 32.119538
$ws.Range("I17").Value = 0.04831535056182164
$ws.Range("J17").Value = 0.032737474561192
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 223.2203288371247
$ws.Range("R17").Value = 1339.321973022748
$ws.Range("S17").Value = 0.00234498918884127
$ws.Range("T17").Value = 0.001636636014182358
$ws.Range("G18").Value = 16.059769
$ws.Range("H18").Value = 32.119538
$ws.Range("I18").Value = 0.04831535056182164
$ws.Range("J18").Value = 0.032737474561192
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 1130.003114063432
$ws.Range("R18").Value = 6780.01868438059
$ws.Range("S18").Value = 0.01187098459911869
$ws.Range("T18").Value = 0.008285104686696639
$ws.Range("G19").Value = 16.059769
$ws.Range("H19").Value = 32.119538
$ws.Range("I19").Value = 0.04831535056182164
$ws.Range("J19").Value = 0.032737474561192
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 1317.024156129385
$ws.Range("R19").Value = 7902.14493677631
$ws.Range("S19").Value = 0.01383569060961153
$ws.Range("T19").Value = 0.009656330033642491
$ws.Range("G20").Value = 16.059769
$ws.Range("H20").Value = 32.119538
$ws.Range("I20").Value = 0.04831535056182164
$ws.Range("J20").Value = 0.032737474561192
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 402.2995662061585
$ws.Range("R20").Value = 1609.198264824634
$ws.Range("S20").Value = 0.004226264419300845
$ws.Range("T20").Value = 0.001966421732205108
$ws.Range("G21").Value = 16.059769
$ws.Range("H21").Value = 32.119538
$ws.Range("I21").Value = 0.04831535056182164
$ws.Range("J21").Value = 0.032737474561192
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 1526.607701494847
$ws.Range("R21").Value = 9159.646208969081
$ws.Range("S21").Value = 0.0160374217449493
$ws.Range("T21").Value = 0.01119298209446541
$ws.Range("G22").Value = 0.5092873333333333
$ws.Range("H22").Value = 1.527862
$ws.Range("I22").Value = 0.001532176212914103
$ws.Range("J22").Value = 0.001557256002810873
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 7.078762217516889
$ws.Range("R22").Value = 63.70885995765201
$ws.Range("S22").Value = 0.0000743642882211112
$ws.Range("T22").Value = 0.00007785149256818966
$ws.Range("G23").Value = 0.5092873333333333
$ws.Range("H23").Value = 1.527862
$ws.Range("I23").Value = 0.001532176212914103
$ws.Range("J23").Value = 0.001557256002810873
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 35.83465444737889
$ws.Range("R23").Value = 322.51189002641
$ws.Range("S23").Value = 0.0003764526183736655
$ws.Range("T23").Value = 0.0003941058123820368
$ws.Range("G24").Value = 0.5092873333333333
$ws.Range("H24").Value = 1.527862
$ws.Range("I24").Value = 0.001532176212914103
$ws.Range("J24").Value = 0.001557256002810873
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 41.76546502074333
$ws.Range("R24").Value = 375.88918518669
$ws.Range("S24").Value = 0.000438757367892035
$ws.Range("T24").Value = 0.0004593322518481144
$ws.Range("G25").Value = 0.5092873333333333
$ws.Range("H25").Value = 1.527862
$ws.Range("I25").Value = 0.001532176212914103
$ws.Range("J25").Value = 0.001557256002810873
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 12.75772231059433
$ws.Range("R25").Value = 76.54633386356601
$ws.Range("S25").Value = 0.0001340232811609729
$ws.Range("T25").Value = 0.00009353873771815647
$ws.Range("G26").Value = 0.5092873333333333
$ws.Range("H26").Value = 1.527862
$ws.Range("I26").Value = 0.001532176212914103
$ws.Range("J26").Value = 0.001557256002810873
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 48.41177761276889
$ws.Range("R26").Value = 435.7059985149201
$ws.Range("S26").Value = 0.0005085786572663182
$ws.Range("T26").Value = 0.0005324277082943753
